# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows into the dataset:
#   - a new row at (original) row 150, shifting the existing rows 150-219 down by one
#   - a new row at (new) row 220, shifting the existing rows 220-230 down by one more
#
# Net effect: dimension grows from A1:T230 to A1:T232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Insert first new row at row 150 ----
$ws.Rows.Item(150).Insert()

$ws.Cells.Item(150, 1).Value = 10
$ws.Cells.Item(150, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(150, 3).Value = "La Araucanía"
$ws.Cells.Item(150, 4).Value = 44567
$ws.Cells.Item(150, 5).Value = 9
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100101
$ws.Cells.Item(150, 8).Value = "Berries"
$ws.Cells.Item(150, 9).Value = 100112025
$ws.Cells.Item(150, 10).Value = "Frutilla"
$ws.Cells.Item(150, 11).Value = "Sin especificar"
$ws.Cells.Item(150, 12).Value = "Primera"
$ws.Cells.Item(150, 13).Value = 130
$ws.Cells.Item(150, 14).Value = 7000
$ws.Cells.Item(150, 15).Value = 8000
$ws.Cells.Item(150, 16).Value = 7538
$ws.Cells.Item(150, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(150, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(150, 19).Value = 1077
$ws.Cells.Item(150, 20).Value = 7

# ---- Insert second new row at (new) row 220 ----
$ws.Rows.Item(220).Insert()

$ws.Cells.Item(220, 1).Value = 10
$ws.Cells.Item(220, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(220, 3).Value = "La Araucanía"
$ws.Cells.Item(220, 4).Value = 44568
$ws.Cells.Item(220, 5).Value = 9
$ws.Cells.Item(220, 6).Value = "Fruta"
$ws.Cells.Item(220, 7).Value = 100101
$ws.Cells.Item(220, 8).Value = "Berries"
$ws.Cells.Item(220, 9).Value = 100112025
$ws.Cells.Item(220, 10).Value = "Frutilla"
$ws.Cells.Item(220, 11).Value = "Sin especificar"
$ws.Cells.Item(220, 12).Value = "Primera"
$ws.Cells.Item(220, 13).Value = 65
$ws.Cells.Item(220, 14).Value = 8000
$ws.Cells.Item(220, 15).Value = 8000
$ws.Cells.Item(220, 16).Value = 8000
$ws.Cells.Item(220, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(220, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(220, 19).Value = 1143
$ws.Cells.Item(220, 20).Value = 7
